$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.789.10"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.01%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.535.42"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.25%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.81"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.11%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "194.56"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.10%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.624"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.26%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.202"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -5.00%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.647"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.94%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.36"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.00%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.23%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.47"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.25%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.089.76"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.75%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "594.19"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -3.34%  "

$ws.Range("B16").Value = "Uniswap"
$ws.Range("C16").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "12.79"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.00%  "

$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.908.75"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.01%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "19.00"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.96%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.536.92"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.22%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.82%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.983"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.55%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "17.63"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.69%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "102.81"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.16%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.10"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.16%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.63"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.17%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.15%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.75"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.84%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.54"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -3.44%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.18"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -3.06%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.04"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.42%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.22"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.09%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.29"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.63%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.114"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.10%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.42"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.26%  "

$ws.Range("B35").Value = "Maker"
$ws.Range("C35").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.820.06"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.10%  "

$ws.Range("B36").Value = "Fetch.AI"
$ws.Range("C36").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.19"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +4.03%  "

$ws.Range("B37").Value = "PEPE"
$ws.Range("C37").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0₃0812"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.57%  "

$ws.Range("B38").Value = "Dai"
$ws.Range("C38").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.18%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "511.35"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.70%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.16%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.56"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.33%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "36.56"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.07%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.133"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.71%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0447"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.38%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.39"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +2.10%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.88%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.85%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.01%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.47"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -3.08%  "

$ws.Range("B50").Value = "FLOKI"
$ws.Range("C50").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000245"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.22%  "

$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.32"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.45%  "
